$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append two new workout entries for Jeremiah (week of 2024-06-22) ---

# Row 93: Jeremiah / Ride
$ws.Cells.Item(93, 1).Value = "Jeremiah"
$ws.Cells.Item(93, 2).Value = 45465
$ws.Cells.Item(93, 3).Value = "Ride"
$ws.Cells.Item(93, 4).Value = 90
$ws.Cells.Item(93, 5).Value = 30.04
$ws.Cells.Item(93, 6).Value = 0
$ws.Cells.Item(93, 7).Value = 7
$ws.Cells.Item(93, 8).Value = 83
$ws.Cells.Item(93, 9).Value = 0
$ws.Cells.Item(93, 10).Value = 0
$ws.Cells.Item(93, 11).Value = 0
$ws.Cells.Item(93, 12).Value = "Agile Antelope"
$ws.Cells.Item(93, 13).Value = 2

# Row 94: Jeremiah / Walk
$ws.Cells.Item(94, 1).Value = "Jeremiah"
$ws.Cells.Item(94, 2).Value = 45465
$ws.Cells.Item(94, 3).Value = "Walk"
$ws.Cells.Item(94, 4).Value = 42
$ws.Cells.Item(94, 5).Value = 1.83
$ws.Cells.Item(94, 6).Value = 171
$ws.Cells.Item(94, 7).Value = 23
$ws.Cells.Item(94, 8).Value = 2
$ws.Cells.Item(94, 9).Value = 0
$ws.Cells.Item(94, 10).Value = 0
$ws.Cells.Item(94, 11).Value = 0
$ws.Cells.Item(94, 12).Value = "Agile Antelope"
$ws.Cells.Item(94, 13).Value = 2

# The Date column (B) uses a short-date number format elsewhere in the
# table (style applied to B2:B92) - copy that format onto the two new
# date cells so they match (rather than Excel inventing a fresh style).
$null = $ws.Range("B92").Copy()
$null = $ws.Range("B93:B94").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Reflect the user's new selection after typing the rows in (Excel
# leaves the cursor one row below the last entry).
$null = $ws.Range("A95").Select()
